$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.932.79'
$ws.Range('E2').Value = '  -0.30%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.044.16'
$ws.Range('E3').Value = '  +0.21%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.67'
$ws.Range('E5').Value = '  -0.26%  '

$ws.Range('E6').Value = '  +0.48%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.31'
$ws.Range('E8').Value = '  -1.18%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.377'
$ws.Range('E9').Value = '  -0.30%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0777'
$ws.Range('E10').Value = '  +2.45%  '

$ws.Range('E11').Value = '  +1.61%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.71'
$ws.Range('E12').Value = '  +4.18%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.340.03'
$ws.Range('E13').Value = '  +0.01%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.59'
$ws.Range('E14').Value = '  +6.56%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.787'
$ws.Range('E15').Value = '  -4.65%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.043.97'
$ws.Range('E16').Value = '  +0.27%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '36.906.94'
$ws.Range('E17').Value = '  -0.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.44'
$ws.Range('E18').Value = '  +13.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.50'
$ws.Range('E19').Value = '  +1.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0888'
$ws.Range('E20').Value = '  +2.99%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.28'
$ws.Range('E21').Value = '  +1.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '234.93'
$ws.Range('E22').Value = '  -1.36%  '

$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.35'
$ws.Range('E24').Value = '  -2.75%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.17'
$ws.Range('E25').Value = '  +8.04%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.93'
$ws.Range('E26').Value = '  -1.76%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.05'
$ws.Range('E27').Value = '  -0.63%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.63'
$ws.Range('E28').Value = '  -3.36%  '

$ws.Range('E29').Value = '  +0.61%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.10'
$ws.Range('E30').Value = '  +1.93%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.65'
$ws.Range('E31').Value = '  +2.18%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0607'
$ws.Range('E32').Value = '  -3.63%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.36'
$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('E34').Value = '  -0.06%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('E36').Value = '  -2.67%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.74'
$ws.Range('E37').Value = '  -1.84%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.106'
$ws.Range('E38').Value = '  -0.82%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.33'
$ws.Range('E39').Value = '  -1.48%  '

$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.20'
$ws.Range('E40').Value = '  +14.01%  '

$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.88'
$ws.Range('E41').Value = '  +22.78%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0219'
$ws.Range('E42').Value = '  -3.09%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.14'
$ws.Range('E43').Value = '  -4.02%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '95.02'
$ws.Range('E44').Value = '  -2.50%  '

$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.10'
$ws.Range('E45').Value = '  -3.31%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.40'
$ws.Range('E46').Value = '  +1.14%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.270.31'
$ws.Range('E47').Value = '  -2.49%  '

$ws.Range('E48').Value = '  -1.46%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.219.98'
$ws.Range('E49').Value = '  -0.32%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.63'
$ws.Range('E50').Value = '  -4.67%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '41.91'
$ws.Range('E51').Value = '  -7.04%  '
